$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.8160071017062
$ws.Range("D2").Value = 2581.94398428306

$ws.Range("B3").Value = 12.8265265956532
$ws.Range("D3").Value = 2252.98315731745

$ws.Range("B4").Value = 12.6149752281368
$ws.Range("D4").Value = 785.833937518773

$ws.Range("B5").Value = 5.91071003380408
$ws.Range("D5").Value = 465.691023207654

$ws.Range("B6").Value = 7.73023881095466
$ws.Range("D6").Value = 424.89906260942

$ws.Range("B7").Value = 2.97744620885537
$ws.Range("D7").Value = 776.260737700715
